# Insert a new column before column AD (EMPLOYEE_ID), shifting
# AD:AJ -> AE:AK. This matches the diff: a new "Unnamed: 0" column
# appears at AD, the old columns shift right by one, and some of the
# shifted values on rows 2-3 are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AD (column 30), pushing existing AD:AJ to AE:AK.
$ws.Columns.Item(30).Insert()

# New column header + blank data cells for rows 2-3.
$ws.Range("AD1").Value = "Unnamed: 0"
$ws.Range("AD2").Value = ""
$ws.Range("AD3").Value = ""

# Update the values that changed (not just shifted) on rows 2 and 3.
$ws.Range("AF2").Value = "O281"
$ws.Range("AI2").Value = "INACTIVE"
$ws.Range("AJ2").Value = "HERO"
$ws.Range("AK2").Value = "RECOVERY"

$ws.Range("AF3").Value = "O72"
$ws.Range("AJ3").Value = "MUTHOOT"
